$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix French/English column order for rows 2-4 (columns A/B and C/D were swapped)
$ws.Range("A2").Value = "comment"
$ws.Range("B2").Value = "how\what"
$ws.Range("C2").Value = "quand"
$ws.Range("D2").Value = "when"

$ws.Range("A3").Value = "quel/quelle"
$ws.Range("B3").Value = "what\how"
$ws.Range("C3").Value = "qui"
$ws.Range("D3").Value = "who"

$ws.Range("A4").Value = "où"
$ws.Range("B4").Value = "where"
$ws.Range("C4").Value = "pourquoi"
$ws.Range("D4").Value = "why"

# Update the active selection to C2:D4
$ws.Range("C2:D4").Select()
